$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Rename field labels (column C) for CompensateAcct and Remark rows, highlighted in red
$ws.Range("C73").Value = "代償專戶"
$ws.Range("C73").Font.Color = 255

$ws.Range("C75").Value = "備註"
$ws.Range("C75").Font.Color = 255

# Mark rows as deleted / renamed in new column H
$ws.Range("H70").Value = "刪除"
$ws.Range("H71").Value = "刪除"
$ws.Range("H72").Value = "刪除"
$ws.Range("H73").Value = "更名"
$ws.Range("H74").Value = "刪除"
$ws.Range("H75").Value = "更名"

# Update the view's top-left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("C68").Select()
